$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($cell, $text) {
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# --- Update the first three data rows ---
Set-CellText $t.Cell(1, 1) "0M"
Set-CellText $t.Cell(2, 1) "0M"
Set-CellText $t.Cell(3, 1) "0M"

# --- Insert 10 new rows before the row currently holding "0" (row 4) ---
$newValues = @("104", "0.00002", "0.00011", "0.00006", "0.00002", "0.00008", "0.00008", "0.00011", "0.00378", "100.0")
$refRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    Set-CellText $newRow.Cells.Item(1) $newValues[$i]
}

# --- Collapse the three trailing multi-run rows into single-run cells ---
$rowCount = $t.Rows.Count
Set-CellText $t.Cell($rowCount - 2, 1) "100"
Set-CellText $t.Cell($rowCount - 1, 1) "0"
Set-CellText $t.Cell($rowCount, 1) "212"
